$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("G3").Value = 3.1
$ws.Range("J3").Value = 3.6
$ws.Range("P3").Value = 2.1
$ws.Range("Q3").Value = 1.6

# Row 4 changes
$ws.Range("F4").Value = 3.15
$ws.Range("G4").Value = 3.85
$ws.Range("I4").Value = 2.48
$ws.Range("J4").Value = 3.55
$ws.Range("K4").Value = 4.9

# Row 5 changes
$ws.Range("F5").Value = 7.4
$ws.Range("G5").Value = 8.199999999999999
$ws.Range("H5").Value = 1.46
$ws.Range("I5").Value = 1.52
$ws.Range("J5").Value = 4.8
$ws.Range("K5").Value = 5.3
$ws.Range("N5").Value = 5.1
$ws.Range("O5").Value = 1.21
$ws.Range("P5").Value = 2.4
$ws.Range("Q5").Value = 1.61
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 1.67
$ws.Range("U5").Value = 2.06
$ws.Range("X5").Value = 980
$ws.Range("Z5").Value = 12.5
$ws.Range("AA5").Value = 15.5
$ws.Range("AB5").Value = 980
$ws.Range("AF5").Value = 70
$ws.Range("AG5").Value = 28
$ws.Range("AH5").Value = 26
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 250
$ws.Range("AK5").Value = 120
$ws.Range("AL5").Value = 85
$ws.Range("AN5").Value = 130
$ws.Range("AO5").Value = 7
